$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '25.934.53'
$ws.Range('E2').Value = '  -0.77%  '
# Row 3
$ws.Range('D3').Value = '1.631.87'
$ws.Range('E3').Value = '  -2.29%  '
# Row 4
$ws.Range('E4').Value = '  -0.02%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.18'
$ws.Range('E5').Value = '  -1.65%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5140'
$ws.Range('E6').Value = '  -1.32%  '
# Row 7
$ws.Range('E7').Value = '  -0.01%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2551'
$ws.Range('E8').Value = '  -3.35%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06187'
$ws.Range('E9').Value = '  -0.90%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.20'
$ws.Range('E10').Value = '  -4.51%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07529'
$ws.Range('E11').Value = '  +0.18%  '
# Row 12
$ws.Range('D12').Value = '1.650.04'
$ws.Range('E12').Value = '  -0.98%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.341'
$ws.Range('E13').Value = '  -2.16%  '
# Row 14
$ws.Range('D14').Value = '1.856.02'
$ws.Range('E14').Value = '  -2.22%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5351'
$ws.Range('E15').Value = '  -4.25%  '
# Row 16
$ws.Range('D16').Value = '0.0₅7914'
$ws.Range('E16').Value = '  -0.79%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.65'
$ws.Range('E17').Value = '  -2.39%  '
# Row 18
$ws.Range('D18').Value = '25.939.34'
$ws.Range('E18').Value = '  -0.92%  '
# Row 19
$ws.Range('E19').Value = '  +0.02%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.608'
$ws.Range('E20').Value = '  -3.86%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '184.81'
$ws.Range('E21').Value = '  -1.16%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.912'
$ws.Range('E22').Value = '  -4.23%  '
# Row 23
$ws.Range('E23').Value = '  +0.03%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.041'
$ws.Range('E24').Value = '  -2.21%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.93'
$ws.Range('E25').Value = '  -1.21%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1189'
$ws.Range('E26').Value = '  -4.38%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.282'
$ws.Range('E27').Value = '  -4.00%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.36'
$ws.Range('E28').Value = '  -3.25%  '
# Row 29
$ws.Range('E29').Value = '  +1.02%  '
# Row 30
$ws.Range('E30').Value = '  -4.41%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.235'
$ws.Range('E31').Value = '  -3.37%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.345'
$ws.Range('E32').Value = '  -3.72%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.324'
$ws.Range('E33').Value = '  -2.95%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.590'
$ws.Range('E34').Value = '  -1.18%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9622'
$ws.Range('E35').Value = '  -2.89%  '
# Row 36
$ws.Range('E36').Value = '  -1.13%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.714'
$ws.Range('E37').Value = '  +0.35%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5768'
$ws.Range('E38').Value = '  -4.31%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01577'
$ws.Range('E39').Value = '  -1.98%  '
# Row 40
$ws.Range('D40').Value = '1.050.39'
$ws.Range('E40').Value = '  -2.04%  '
# Row 41
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.754'
$ws.Range('E41').Value = '  -6.07%  '
# Row 42
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.002'
$ws.Range('E42').Value = '  -0.27%  '
# Row 43
$ws.Range('E43').Value = '  -3.55%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.65'
$ws.Range('E44').Value = '  +0.12%  '
# Row 45
$ws.Range('D45').Value = '1.783.98'
$ws.Range('E45').Value = '  -1.95%  '
# Row 46
$ws.Range('D46').Value = '0.0₈107'
$ws.Range('E46').Value = '  -0.44%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.004'
$ws.Range('E47').Value = '  +0.09%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.89'
$ws.Range('E48').Value = '  -3.68%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05195'
$ws.Range('E49').Value = '  -1.07%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.869'
$ws.Range('E50').Value = '  -0.94%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4221'
$ws.Range('E51').Value = '  -0.72%  '
